# Sync attendance_reports: normalise the "Recorded By" (column G) list
# ordering. Multi-author cells ("A, B" / "A, B, C") have their
# comma-separated entries reversed in place; single-author cells are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Text

    if ($val -ne $null -and $val -match ",") {
        $parts = $val -split ", "
        $n = $parts.Count

        $reversed = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }

        $newVal = $reversed -join ", "
        $cell.Value = $newVal
    }
}
